$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.168.94"
$ws.Range("E2").Value = "  -3.23%  "

# Row 3
$ws.Range("D3").Value = "1.715.21"
$ws.Range("E3").Value = "  -3.60%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.56"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.03%  "

# Row 6
$ws.Range("E6").Value = "  +0.04%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4789"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +6.46%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3456"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.12"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07275"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.37%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.046"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.87"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.863"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.35%  "

# Row 15
$ws.Range("D15").Value = "1.713.87"
$ws.Range("E15").Value = "  -3.56%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.872"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.94"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.27%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001040"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.40%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06366"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.23%  "

# Row 20
$ws.Range("E20").Value = "  +0.06%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.49"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.91%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.613"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.25%  "

# Row 23
$ws.Range("D23").Value = "27.199.45"
$ws.Range("E23").Value = "  -3.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.83"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.50%  "

# Row 25
$ws.Range("E25").Value = "  -1.62%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.15"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.64"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.80%  "

# Row 28
$ws.Range("D28").Value = "1.909.18"
$ws.Range("E28").Value = "  -3.65%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.101"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.02%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.99"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.17%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.019"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -8.10%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09250"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.33%  "

# Row 33
$ws.Range("E33").Value = "  -3.01%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.324"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.97%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02199"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.25%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05918"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.75%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.09"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.94%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2006"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.19%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.417"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.41%  "

# Row 40
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.750"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.24%  "

# Row 41
$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9997"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5938"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.33%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.093"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -8.05%  "

# Row 44
$ws.Range("E44").Value = "  -5.82%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.61"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.96%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.576"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.76%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5619"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.04%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.72"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.42%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.840"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -6.37%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06643"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.085"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.17%  "
